# Auto-generated Excel COM-interop script to apply the Ultros_Profits.xlsx cell updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 993.86664
$ws.Range("I33").Value = 695.9
$ws.Range("K33").Value = 695.9
$ws.Range("M33").Value = -466.9

$ws.Range("H62").Value = 2205.3333
$ws.Range("I62").Value = 1856
$ws.Range("K62").Value = 1856
$ws.Range("M62").Value = -1232

$ws.Range("H65").Value = 2205.3333
$ws.Range("I65").Value = 1856
$ws.Range("K65").Value = 9280
$ws.Range("M65").Value = -6160

$ws.Range("H70").Value = 2097.1904
$ws.Range("I70").Value = 1489
$ws.Range("J70").Value = 2766.2
$ws.Range("K70").Value = 4467
$ws.Range("L70").Value = 8298.599999999999
$ws.Range("M70").Value = -4197
$ws.Range("N70").Value = -8838.599999999999

$ws.Range("H73").Value = 2097.1904
$ws.Range("I73").Value = 1489
$ws.Range("J73").Value = 2766.2
$ws.Range("K73").Value = 4467
$ws.Range("L73").Value = 8298.599999999999
$ws.Range("M73").Value = -3531
$ws.Range("N73").Value = -10170.6

$ws.Range("H137").Value = 3322.95
$ws.Range("I137").Value = 3129.4211
$ws.Range("J137").Value = 7000
$ws.Range("K137").Value = 9388.263300000001
$ws.Range("L137").Value = 21000
$ws.Range("M137").Value = -6838.263300000001
$ws.Range("N137").Value = -26100

$ws.Range("H138").Value = 3746.8147
$ws.Range("J138").Value = 5341.2334
$ws.Range("L138").Value = 16023.7002
$ws.Range("N138").Value = -26303.7002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2829.9375
$ws.Range("I74").Value = 2519.9285
$ws.Range("K74").Value = 2519.9285
$ws.Range("M74").Value = -1645.9285

$ws.Range("H77").Value = 2829.9375
$ws.Range("I77").Value = 2519.9285
$ws.Range("K77").Value = 12599.6425
$ws.Range("M77").Value = -8231.6425

$ws.Range("H102").Value = 9701.637000000001
$ws.Range("I102").Value = 9570.700000000001
$ws.Range("K102").Value = 9570.700000000001
$ws.Range("M102").Value = -7948.700000000001

$ws.Range("H122").Value = 2203.8
$ws.Range("I122").Value = 2004.0714
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 6012.2142
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -3562.2142
$ws.Range("N122").Value = -19900

$ws.Range("H132").Value = 1514.8182
$ws.Range("I132").Value = 1514.8182
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4544.4546
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2014.4546
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1196.6578
$ws.Range("I134").Value = 984.64703
$ws.Range("K134").Value = 2953.94109
$ws.Range("M134").Value = -418.9410899999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1791.7307
$ws.Range("I31").Value = 1791.7307
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1791.7307
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1496.7307
$ws.Range("N31").ClearContents()

$ws.Range("H34").Value = 1791.7307
$ws.Range("I34").Value = 1791.7307
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1791.7307
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1589.7307
$ws.Range("N34").ClearContents()

$ws.Range("H51").Value = 23333.334
$ws.Range("J51").Value = 49000
$ws.Range("L51").Value = 49000
$ws.Range("N51").Value = -50472

$ws.Range("H58").Value = 2203.5454
$ws.Range("I58").Value = 1727.2222
$ws.Range("J58").Value = 4347
$ws.Range("K58").Value = 1727.2222
$ws.Range("L58").Value = 4347
$ws.Range("M58").Value = -1524.2222
$ws.Range("N58").Value = -4753

$ws.Range("H61").Value = 23333.334
$ws.Range("J61").Value = 49000
$ws.Range("L61").Value = 49000
$ws.Range("N61").Value = -49696

$ws.Range("H86").Value = 17298.152
$ws.Range("I86").Value = 32870.637
$ws.Range("J86").Value = 9511.909
$ws.Range("K86").Value = 32870.637
$ws.Range("L86").Value = 9511.909
$ws.Range("M86").Value = -31747.637
$ws.Range("N86").Value = -11757.909

$ws.Range("H89").Value = 17298.152
$ws.Range("I89").Value = 32870.637
$ws.Range("J89").Value = 9511.909
$ws.Range("K89").Value = 164353.185
$ws.Range("L89").Value = 47559.545
$ws.Range("M89").Value = -158737.185
$ws.Range("N89").Value = -58791.545

$ws.Range("H94").Value = 14088.75
$ws.Range("I94").Value = 34041.332
$ws.Range("J94").Value = 2117.2
$ws.Range("K94").Value = 34041.332
$ws.Range("L94").Value = 2117.2
$ws.Range("M94").Value = -33590.332
$ws.Range("N94").Value = -3019.2

$ws.Range("H99").Value = 56100260
$ws.Range("I99").Value = 8133681.5
$ws.Range("K99").Value = 8133681.5
$ws.Range("M99").Value = -8132183.5

$ws.Range("H126").Value = 56100260
$ws.Range("I126").Value = 8133681.5
$ws.Range("K126").Value = 24401044.5
$ws.Range("M126").Value = -24398574.5

$ws.Range("H132").Value = 3001.125
$ws.Range("I132").Value = 2787
$ws.Range("K132").Value = 8361
$ws.Range("M132").Value = -5831

$ws.Range("H134").Value = 3773.7827
$ws.Range("I134").Value = 2608.0952
$ws.Range("K134").Value = 7824.285600000001
$ws.Range("M134").Value = -5289.285600000001

$ws.Range("H136").Value = 2203.5454
$ws.Range("I136").Value = 1727.2222
$ws.Range("J136").Value = 4347
$ws.Range("K136").Value = 5181.6666
$ws.Range("L136").Value = 13041
$ws.Range("M136").Value = -2631.6666
$ws.Range("N136").Value = -18141

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 53.25
$ws.Range("I17").Value = 50.5
$ws.Range("J17").Value = 56
$ws.Range("K17").Value = 151.5
$ws.Range("L17").Value = 168
$ws.Range("M17").Value = 17.5
$ws.Range("N17").Value = -506

$ws.Range("H95").Value = 5000
$ws.Range("J95").Value = 5000
$ws.Range("L95").Value = 15000
$ws.Range("N95").Value = -19118

$ws.Range("H113").Value = 1287.9231
$ws.Range("J113").Value = 1221.7778
$ws.Range("L113").Value = 3665.3334
$ws.Range("N113").Value = -8005.3334

$ws.Range("H132").Value = 1763.909
$ws.Range("I132").Value = 1555.8889
$ws.Range("K132").Value = 14003.0001
$ws.Range("M132").Value = -11473.0001

$ws.Range("H137").Value = 2657.9546
$ws.Range("I137").Value = 2305.7693
$ws.Range("K137").Value = 6917.3079
$ws.Range("M137").Value = -1817.3079

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2157.2163
$ws.Range("I132").Value = 2495.6843
$ws.Range("J132").Value = 1799.9445
$ws.Range("K132").Value = 7487.0529
$ws.Range("L132").Value = 5399.833500000001
$ws.Range("M132").Value = -4957.0529
$ws.Range("N132").Value = -10459.8335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 12857
$ws.Range("I40").Value = 14977
$ws.Range("K40").Value = 14977
$ws.Range("M40").Value = -14841

$ws.Range("H93").Value = 8140.75
$ws.Range("J93").Value = 9432.799999999999
$ws.Range("L93").Value = 9432.799999999999
$ws.Range("N93").Value = -11928.8

$ws.Range("H100").Value = 106118.63
$ws.Range("I100").Value = 188634.17
$ws.Range("K100").Value = 188634.17
$ws.Range("M100").Value = -188093.17

$ws.Range("H122").Value = 7681.8184
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws.Range("H136").Value = 3737.8125
$ws.Range("I136").Value = 4099.615
$ws.Range("J136").Value = 2170
$ws.Range("K136").Value = 12298.845
$ws.Range("L136").Value = 6510
$ws.Range("M136").Value = -9748.844999999999
$ws.Range("N136").Value = -11610

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1098.2858
$ws.Range("I107").Value = 1030
$ws.Range("J107").Value = 1149.5
$ws.Range("K107").Value = 3090
$ws.Range("L107").Value = 3448.5
$ws.Range("M107").Value = -1170
$ws.Range("N107").Value = -7288.5

$ws.Range("H122").Value = 3333.3333

$ws.Range("H132").Value = 2487.6
$ws.Range("I132").Value = 2568.4285
$ws.Range("J132").Value = 2299
$ws.Range("K132").Value = 7705.2855
$ws.Range("L132").Value = 6897
$ws.Range("M132").Value = -5175.2855
$ws.Range("N132").Value = -11957
